$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3

$ws.Range("E5").Select()
